# Restore the "R20" rule's lower bound (column C, row 10) from 18 to 20
# on the "Rules" worksheet, matching the committed revision.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Cells.Item(10, 3).Value = 20
